$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add photo_url values in column B for rows 2-8
$ws.Range("B2").Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-debbie-gray.jpg?h=1200&la=en&w=1200&rev=5cf7cf91a33849419def2bc4d4d030f3"
$ws.Range("B3").Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-harpreet-kaur.jpg?h=1200&la=en&w=1200&rev=347fbb20f8d34085abb16563707b2cde"
$ws.Range("B4").Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-madison-mcguire.jpg?h=1200&la=en&w=1200&rev=cffce2fefade4171ae6a7fa8e3b9015c"
$ws.Range("B5").Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-savanah-nunes.jpg?h=1200&la=en&w=1200&rev=dfa96124f83244bb910bcdddfc069582"
$ws.Range("B6").Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-brandi-chastain.jpg?h=1200&la=en&w=1200&rev=22142ccc9d92424bb2a8a749cfa522ff"
$ws.Range("B7").Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-brian-kelly.jpg?h=300&la=en&w=300&rev=7a40f27403394619b8297896d9959f5c"
$ws.Range("B8").Value = "/-/media/uidaho-responsive/images/cals/college/our-people/1200x1200-carly-schoepflin.jpg?h=1200&la=en&w=1200&rev=b17fed7c84c84bd3bd657141bd56205f"

# Clear title/office/phone/email for row 4 (Madison McGuire)
$ws.Range("F4:I4").ClearContents()

# Clear title/office/phone/email for row 8 (Carly Schoepflin)
$ws.Range("F8:I8").ClearContents()
